$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table runs through row 43 (period "01-06-2021"); append the next
# month ("01-07-2021") as a new row 44, same 13-column layout (A = period
# label, B:M = the 12 division index values).
$row = 44

# Column A stores the period as literal text (e.g. "01-06-2021"), like every
# other row above it. A plain $cell.Value = "01-07-2021" assignment makes
# Excel "helpfully" interpret the text as a date literal (storing a date
# serial number plus a new date number-format style) instead of keeping it
# as text. To avoid that, write it as a formula that evaluates to the
# literal string, then Copy / Paste Special-Values it back onto itself -
# this flattens the formula to a plain static text cell (matching the
# existing rows) without touching any cell styles.
$cell = $ws.Cells.Item($row, 1)
$cell.Formula = '="01-07-2021"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item($row, 2).Value = 114.7
$ws.Cells.Item($row, 3).Value = 109.27
$ws.Cells.Item($row, 4).Value = 98.47
$ws.Cells.Item($row, 5).Value = 111.1
$ws.Cells.Item($row, 6).Value = 111.5
$ws.Cells.Item($row, 7).Value = 107.27
$ws.Cells.Item($row, 8).Value = 110.74
$ws.Cells.Item($row, 9).Value = 93.41
$ws.Cells.Item($row, 10).Value = 109.05
$ws.Cells.Item($row, 11).Value = 112.12
$ws.Cells.Item($row, 12).Value = 108.96
$ws.Cells.Item($row, 13).Value = 110.68
